$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

# --- Row 44 / 45: BitcoinSV and FirstDigitalUSD swap ranking positions, with refreshed price/volume ---
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D44") "1.00"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D45") "93.75"
$ws.Range("E45").Value = "  +46.54%  "

# --- Price (D) / Volume 1h change (E) refresh for the remaining rows ---
$ws.Range("D2").Value = "42.956.30"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.400.62"
$ws.Range("E3").Value = "  +4.77%  "
$ws.Range("E4").Value = "  -0.30%  "
Set-TextValue $ws.Range("D5") "335.71"
$ws.Range("E5").Value = "  +8.45%  "
Set-TextValue $ws.Range("D6") "100.68"
$ws.Range("E6").Value = "  -10.48%  "
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("E8").Value = "  -0.04%  "
Set-TextValue $ws.Range("D9") "0.638"
$ws.Range("E9").Value = "  +3.51%  "
Set-TextValue $ws.Range("D10") "40.58"
$ws.Range("E10").Value = "  -8.69%  "
$ws.Range("E11").Value = "  +0.38%  "
Set-TextValue $ws.Range("D12") "8.53"
$ws.Range("E12").Value = "  -3.48%  "
$ws.Range("E13").Value = "  -4.38%  "
Set-TextValue $ws.Range("D14") "16.93"
$ws.Range("E14").Value = "  +8.94%  "
$ws.Range("E15").Value = "  +1.73%  "
$ws.Range("D16").Value = "2.759.84"
$ws.Range("E16").Value = "  +4.64%  "
$ws.Range("D17").Value = "2.400.93"
$ws.Range("E17").Value = "  +4.97%  "
$ws.Range("D18").Value = "42.886.38"
$ws.Range("E18").Value = "  -0.02%  "
Set-TextValue $ws.Range("D19") "7.70"
$ws.Range("E19").Value = "  +6.75%  "
$ws.Range("E20").Value = "  -0.83%  "
Set-TextValue $ws.Range("D21") "3.89"
$ws.Range("E21").Value = "  +10.12%  "
Set-TextValue $ws.Range("D22") "76.45"
$ws.Range("E22").Value = "  +0.20%  "
Set-TextValue $ws.Range("D23") "269.71"
$ws.Range("E23").Value = "  +4.95%  "
Set-TextValue $ws.Range("D24") "2.37"
$ws.Range("E24").Value = "  -3.73%  "
Set-TextValue $ws.Range("D25") "10.18"
$ws.Range("E25").Value = "  +12.65%  "
Set-TextValue $ws.Range("D26") "11.79"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.05%  "
Set-TextValue $ws.Range("D28") "24.28"
$ws.Range("E28").Value = "  +8.50%  "
Set-TextValue $ws.Range("D30") "174.28"
$ws.Range("E30").Value = "  +0.45%  "
Set-TextValue $ws.Range("D31") "3.11"
$ws.Range("E31").Value = "  -2.05%  "
Set-TextValue $ws.Range("D32") "0.0921"
$ws.Range("E32").Value = "  +2.14%  "
Set-TextValue $ws.Range("D33") "35.80"
$ws.Range("E33").Value = "  -8.51%  "
Set-TextValue $ws.Range("D34") "5.99"
$ws.Range("E34").Value = "  +4.08%  "
$ws.Range("E35").Value = "  +3.16%  "
Set-TextValue $ws.Range("D36") "4.76"
$ws.Range("E36").Value = "  -6.71%  "
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("E38").Value = "  -5.79%  "
Set-TextValue $ws.Range("D39") "0.107"
$ws.Range("E39").Value = "  +2.14%  "
$ws.Range("E40").Value = "  +11.01%  "
$ws.Range("E41").Value = "  +6.81%  "
$ws.Range("E42").Value = "  +0.36%  "
Set-TextValue $ws.Range("D43") "69.60"
$ws.Range("E43").Value = "  -3.54%  "
Set-TextValue $ws.Range("D46") "118.63"
$ws.Range("E46").Value = "  +9.37%  "
Set-TextValue $ws.Range("D47") "11.96"
$ws.Range("E47").Value = "  -3.59%  "
$ws.Range("E48").Value = "  -3.45%  "
Set-TextValue $ws.Range("D49") "9.13"
$ws.Range("E49").Value = "  +2.25%  "
$ws.Range("D50").Value = "1.643.66"
$ws.Range("E50").Value = "  +11.07%  "
Set-TextValue $ws.Range("D51") "1.28"
$ws.Range("E51").Value = "  -2.04%  "
